# Update New Orleans xlsx file:
#  1. Insert a new "State" column into the hotel_info sheet, right after
#     "Hotel_Name" and before "City", populated with "Louisiana" for the
#     existing hotel row.
#  2. Reorder the worksheet tabs so "review_info" comes before "hotel_info".

$wb = $excel.ActiveWorkbook

$hotelInfo  = $wb.Worksheets.Item("hotel_info")
$reviewInfo = $wb.Worksheets.Item("review_info")

# --- 1. Add the new "State" column (column C) to hotel_info ---------------
$hotelInfo.Columns.Item(3).Insert()
$hotelInfo.Cells.Item(1, 3).Value = "State"
$hotelInfo.Cells.Item(2, 3).Value = "Louisiana"

# --- 2. Move hotel_info so it follows review_info --------------------------
$hotelInfo.Move($null, $reviewInfo)
